$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated coin data
$ws.Range("D2").Value = "42.633.18"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "2.356.41"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "319.51"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "108.42"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("D10").Value = "41.64"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "0.997"
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "16.01"
$ws.Range("E15").Value = "  -5.30%  "
$ws.Range("D16").Value = "2.715.35"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "2.364.14"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "42.654.70"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "76.41"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "3.59"
$ws.Range("E22").Value = "  +6.40%  "
$ws.Range("D23").Value = "257.15"
$ws.Range("E23").Value = "  -7.07%  "
$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").Value = "9.42"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "22.84"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("D30").Value = "174.80"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "36.87"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("E33").Value = "  +4.29%  "
$ws.Range("D34").Value = "2.86"
$ws.Range("E34").Value = "  -10.34%  "
$ws.Range("E35").Value = "  +19.96%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  -4.74%  "
$ws.Range("D38").Value = "0.0364"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "3.86"
$ws.Range("E39").Value = "  -7.02%  "
$ws.Range("D40").Value = "2.68"
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.238"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "71.46"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("E43").Value = "  -6.40%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "12.03"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "113.16"
$ws.Range("E46").Value = "  -7.52%  "
$ws.Range("D47").Value = "5.50"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D49").Value = "84.83"
$ws.Range("E49").Value = "  -10.85%  "
$ws.Range("D50").Value = "73.88"
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("D51").Value = "1.28"
$ws.Range("E51").Value = "  -1.33%  "
